# Remove the "Architecture" diagram slide (the second slide in the deck,
# containing the UI / Logic / Storage / Model / Commons rounded-rectangle
# diagram) from the Diagrams.pptx deck.
$p = $ppt.ActivePresentation

$s = $p.Slides.Item(2)
$s.Delete()
